# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update dSF (column F) values per repulled data
$ws.Range("F2").Value = -9
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -5
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -2
$ws.Range("F14").Value = -8
